$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Checklist updates: "command factory pattern" column (E) and
# --- "undo/redo" column (C) for "2. Add hero to player" (row 4) and
# --- "7. delete hero" (row 9) are now implemented ("OK").
$ws.Range("C4").Value = "Y (OK)"
$ws.Range("E4").Value = "OK"
$ws.Range("E9").Value = "OK"

# --- Highlight the two exceptions not covered by the command factory
# --- work: "6. call hero skill" (row 8) and
# --- "8. change name of current player" (row 10).
$ws.Range("A8").Interior.Color = 65535
$ws.Range("A10").Interior.Color = 65535

# --- Freeze panes at column C and move the selection to C9.
$ws.Range("C1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C9").Select()
